# Update FFXIV Leve profit calculations (H-N columns) across all job sheets
# Values refreshed from current Universalis market data by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2483.1667
$ws.Range("J32").Value = 2445.7273
$ws.Range("L32").Value = 2445.7273
$ws.Range("N32").Value = -3097.7273

$ws.Range("H103").Value = 800
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2400
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -3572

$ws.Range("H137").Value = 20839262
$ws.Range("I137").Value = 31257654
$ws.Range("K137").Value = 93772962
$ws.Range("M137").Value = -93770412

$ws.Range("H138").Value = 3340.79
$ws.Range("J138").Value = 4351.7188
$ws.Range("L138").Value = 13055.1564
$ws.Range("N138").Value = -23335.1564

$ws.Range("H141").Value = 5917.0835
$ws.Range("I141").Value = 3213.275
$ws.Range("K141").Value = 9639.825000000001
$ws.Range("M141").Value = -4459.825000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 492.4091
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 458.25
$ws.Range("K97").Value = 500
$ws.Range("L97").Value = 458.25
$ws.Range("M97").Value = -4
$ws.Range("N97").Value = -1450.25

$ws.Range("H132").Value = 14176.027
$ws.Range("I132").Value = 9044.767
$ws.Range("J132").Value = 39832.332
$ws.Range("K132").Value = 27134.301
$ws.Range("L132").Value = 119496.996
$ws.Range("M132").Value = -24604.301
$ws.Range("N132").Value = -124556.996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2430.36
$ws.Range("I86").Value = 2119.6758
$ws.Range("K86").Value = 2119.6758
$ws.Range("M86").Value = -996.6758

$ws.Range("H89").Value = 2430.36
$ws.Range("I89").Value = 2119.6758
$ws.Range("K89").Value = 10598.379
$ws.Range("M89").Value = -4982.379000000001

$ws.Range("H94").Value = 2323.9268
$ws.Range("I94").Value = 1295.5186
$ws.Range("K94").Value = 1295.5186
$ws.Range("M94").Value = -844.5186000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1735.0869
$ws.Range("I16").Value = 1136.3572
$ws.Range("J16").Value = 2666.4443
$ws.Range("K16").Value = 1136.3572
$ws.Range("L16").Value = 2666.4443
$ws.Range("M16").Value = -849.3571999999999
$ws.Range("N16").Value = -3240.4443

$ws.Range("H31").Value = 273795.53
$ws.Range("I31").Value = 4535.439
$ws.Range("K31").Value = 4535.439
$ws.Range("M31").Value = -4240.439

$ws.Range("H34").Value = 273795.53
$ws.Range("I34").Value = 4535.439
$ws.Range("K34").Value = 4535.439
$ws.Range("M34").Value = -4333.439

$ws.Range("H58").Value = 2156.4814
$ws.Range("I58").Value = 2265.1428
$ws.Range("J58").Value = 2039.4615
$ws.Range("K58").Value = 2265.1428
$ws.Range("L58").Value = 2039.4615
$ws.Range("M58").Value = -2062.1428
$ws.Range("N58").Value = -2445.4615

$ws.Range("H94").Value = 920.4286
$ws.Range("I94").Value = 933.1429000000001
$ws.Range("J94").Value = 914.0714
$ws.Range("K94").Value = 933.1429000000001
$ws.Range("L94").Value = 914.0714
$ws.Range("M94").Value = -482.1429000000001
$ws.Range("N94").Value = -1816.0714

$ws.Range("H113").Value = 1735.0869
$ws.Range("I113").Value = 1136.3572
$ws.Range("J113").Value = 2666.4443
$ws.Range("K113").Value = 1136.3572
$ws.Range("L113").Value = 2666.4443
$ws.Range("M113").Value = 1033.6428
$ws.Range("N113").Value = -7006.4443

$ws.Range("H136").Value = 2156.4814
$ws.Range("I136").Value = 2265.1428
$ws.Range("J136").Value = 2039.4615
$ws.Range("K136").Value = 6795.428400000001
$ws.Range("L136").Value = 6118.3845
$ws.Range("M136").Value = -4245.428400000001
$ws.Range("N136").Value = -11218.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4749.5
$ws.Range("J55").Value = 4999.5
$ws.Range("L55").Value = 14998.5
$ws.Range("N55").Value = -15352.5

$ws.Range("H122").Value = 1172.9445
$ws.Range("I122").Value = 2173.8333
$ws.Range("J122").Value = 672.5
$ws.Range("K122").Value = 19564.4997
$ws.Range("L122").Value = 6052.5
$ws.Range("M122").Value = -17114.4997
$ws.Range("N122").Value = -10952.5

$ws.Range("H139").Value = 16613.5
$ws.Range("I139").Value = 17873.95
$ws.Range("J139").Value = 13462.375
$ws.Range("K139").Value = 53621.85000000001
$ws.Range("L139").Value = 40387.125
$ws.Range("M139").Value = -48481.85000000001
$ws.Range("N139").Value = -50667.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 19000
$ws.Range("J55").Value = 27000
$ws.Range("L55").Value = 27000
$ws.Range("N55").Value = -27654

$ws.Range("H70").Value = 20433.549
$ws.Range("I70").Value = 28437.191
$ws.Range("K70").Value = 28437.191
$ws.Range("M70").Value = -28167.191

$ws.Range("H73").Value = 20433.549
$ws.Range("I73").Value = 28437.191
$ws.Range("K73").Value = 28437.191
$ws.Range("M73").Value = -27501.191

$ws.Range("H107").Value = 2255.875
$ws.Range("I107").Value = 2006.7142
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 2006.7142
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -86.71419999999989
$ws.Range("N107").Value = -7840

$ws.Range("H132").Value = 4636.1665
$ws.Range("I132").Value = 3738.7942
$ws.Range("J132").Value = 8450
$ws.Range("K132").Value = 11216.3826
$ws.Range("L132").Value = 25350
$ws.Range("M132").Value = -8686.382599999999
$ws.Range("N132").Value = -30410

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 631.6
$ws.Range("I46").Value = 688
$ws.Range("J46").Value = 617.5
$ws.Range("K46").Value = 688
$ws.Range("L46").Value = 617.5
$ws.Range("M46").Value = -500
$ws.Range("N46").Value = -993.5

$ws.Range("H93").Value = 2617.5881
$ws.Range("I93").Value = 1212.2667
$ws.Range("J93").Value = 13157.5
$ws.Range("K93").Value = 1212.2667
$ws.Range("L93").Value = 13157.5
$ws.Range("M93").Value = 35.7333000000001
$ws.Range("N93").Value = -15653.5

$ws.Range("H136").Value = 5822.1304
$ws.Range("I136").Value = 5614.7144
$ws.Range("K136").Value = 16844.1432
$ws.Range("M136").Value = -14294.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 11495388
$ws.Range("I126").Value = 16667883
$ws.Range("K126").Value = 50003649
$ws.Range("M126").Value = -50001179

$ws.Range("H132").Value = 1622
$ws.Range("I132").Value = 1433.6522
$ws.Range("K132").Value = 4300.9566
$ws.Range("M132").Value = -1770.9566

$ws.Range("H139").Value = 45399.8
$ws.Range("I139").Value = 29666.334
$ws.Range("K139").Value = 29666.334
$ws.Range("M139").Value = -24526.334
